$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.288.78'
$ws.Range('E2').Value = '  +5.90%  '
$ws.Range('D3').Value = '2.621.77'
$ws.Range('E3').Value = '  +9.42%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '507.82'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.25'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.96%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.995'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('E8').Value = '  -0.82%  '
$ws.Range('D9').Value = '2.662.08'
$ws.Range('E9').Value = '  +10.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.41'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.47%  '
$ws.Range('E11').Value = '  +5.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.346'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.77%  '
$ws.Range('E13').Value = '  +1.10%  '
$ws.Range('D14').Value = '3.104.82'
$ws.Range('E14').Value = '  +10.42%  '
$ws.Range('D15').Value = '60.431.72'
$ws.Range('E15').Value = '  +6.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.81'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.21%  '
$ws.Range('E17').Value = '  +6.27%  '
$ws.Range('D18').Value = '2.652.02'
$ws.Range('E18').Value = '  +10.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.80'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '349.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +7.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.52'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.19'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.97%  '
$ws.Range('E25').Value = '  +5.28%  '
$ws.Range('D26').Value = '2.749.38'
$ws.Range('E26').Value = '  +9.66%  '
$ws.Range('E27').Value = '  +4.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.991'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.61%  '
$ws.Range('E29').Value = '  +12.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.54'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.55%  '
$ws.Range('E31').Value = '  -0.11%  '
$ws.Range('E32').Value = '  +5.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '156.66'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.85%  '
$ws.Range('E34').Value = '  +3.56%  '
$ws.Range('E35').Value = '  +9.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.06'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.21'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '309.08'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +15.45%  '
$ws.Range('E39').Value = '  +9.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.855'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.842'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +32.58%  '
$ws.Range('E42').Value = '  +7.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '35.22'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.641'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +8.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0580'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +9.93%  '
$ws.Range('E46').Value = '  -0.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '20.19'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +16.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.992'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.90'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +9.46%  '
$ws.Range('D50').Value = '2.057.55'
$ws.Range('E50').Value = '  +10.38%  '
$ws.Range('E51').Value = '  +3.55%  '
